$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A16").Value = "r775"
$ws.Range("B16").Value = "test"
$ws.Range("C16").Value = "quick second test"
$ws.Range("D16").Value = "2025-10-01 14:55:11"
